$d = $word.ActiveDocument

# Paragraph styles "Normal" and "Heading" switch their east-Asian font from
# "DejaVu Sans" to "Tahoma" (w:rPr/w:rFonts/@w:eastAsia). That attribute is
# reachable on a Style object via Font.NameFarEast.

$normal = $d.Styles("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $d.Styles("Heading")
$heading.Font.NameFarEast = "Tahoma"

# Paragraph styles "List", "Caption" and "Index" gain an explicit
# complex-script font (w:rPr/w:rFonts/@w:cs = "DejaVu Sans"), which were
# previously empty (<w:rPr/>) or missing an rFonts child. That attribute is
# reachable on a Style object via Font.NameBi.

$list = $d.Styles("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $d.Styles("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $d.Styles("Index")
$index.Font.NameBi = "DejaVu Sans"
